# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# Values are prefixed with a leading apostrophe so Excel stores them
# as literal text (preserves trailing zeros / dotted "thousands" look,
# and the padded spacing already present in the Volume(1h) column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.763.07"
$ws.Range("E2").Value = "'  +0.25%  "
$ws.Range("D3").Value = "'1.636.22"
$ws.Range("E3").Value = "'  -0.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'212.22"
$ws.Range("E5").Value = "'  -0.23%  "
$ws.Range("E6").Value = "'  -0.22%  "
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("D8").Value = "'23.41"
$ws.Range("E8").Value = "'  +1.07%  "
$ws.Range("E9").Value = "'  +2.15%  "
$ws.Range("D10").Value = "'0.0612"
$ws.Range("E10").Value = "'  +0.16%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "'  -3.24%  "
$ws.Range("D12").Value = "'1.868.44"
$ws.Range("E12").Value = "'  -0.11%  "
$ws.Range("D13").Value = "'1.636.99"
$ws.Range("E13").Value = "'  -0.18%  "
$ws.Range("E14").Value = "'  -0.20%  "
$ws.Range("E15").Value = "'  -1.26%  "
$ws.Range("D16").Value = "'65.17"
$ws.Range("E16").Value = "'  +0.68%  "
$ws.Range("D17").Value = "'27.720.89"
$ws.Range("E17").Value = "'  +0.22%  "
$ws.Range("D18").Value = "'231.07"
$ws.Range("E18").Value = "'  +0.24%  "
$ws.Range("E19").Value = "'  -0.33%  "
$ws.Range("E20").Value = "'  -1.29%  "
$ws.Range("E21").Value = "'  -0.07%  "
$ws.Range("D22").Value = "'10.73"
$ws.Range("E22").Value = "'  +4.56%  "
$ws.Range("E23").Value = "'  +1.19%  "
$ws.Range("E24").Value = "'  +2.62%  "
$ws.Range("D25").Value = "'148.91"
$ws.Range("E25").Value = "'  -1.21%  "
$ws.Range("E26").Value = "'  -0.73%  "
$ws.Range("E27").Value = "'  -0.78%  "
$ws.Range("D28").Value = "'15.58"
$ws.Range("E28").Value = "'  -0.07%  "
$ws.Range("E29").Value = "'  -0.03%  "
$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("E31").Value = "'  -0.53%  "
$ws.Range("D32").Value = "'3.29"
$ws.Range("D33").Value = "'1.480.99"
$ws.Range("E33").Value = "'  +1.48%  "
$ws.Range("D34").Value = "'3.10"
$ws.Range("E34").Value = "'  -0.79%  "
$ws.Range("E35").Value = "'  -0.99%  "
$ws.Range("E36").Value = "'  -1.80%  "
$ws.Range("D37").Value = "'0.955"
$ws.Range("E37").Value = "'  +6.77%  "
$ws.Range("E38").Value = "'  +0.31%  "
$ws.Range("E39").Value = "'  -1.29%  "
$ws.Range("E40").Value = "'  +0.11%  "
$ws.Range("D42").Value = "'68.01"
$ws.Range("E42").Value = "'  -1.82%  "
$ws.Range("D43").Value = "'2.45"
$ws.Range("E43").Value = "'  -0.25%  "
$ws.Range("D44").Value = "'2.20"
$ws.Range("E44").Value = "'  -1.34%  "
$ws.Range("E45").Value = "'  -4.31%  "
$ws.Range("D46").Value = "'1.776.84"
$ws.Range("E46").Value = "'  -0.20%  "
$ws.Range("E47").Value = "'  +1.06%  "
$ws.Range("D48").Value = "'87.76"
$ws.Range("E48").Value = "'  +1.13%  "
$ws.Range("E49").Value = "'  -1.66%  "
$ws.Range("D50").Value = "'0.0993"
$ws.Range("E50").Value = "'  +0.06%  "
$ws.Range("D51").Value = "'7.74"
$ws.Range("E51").Value = "'  -0.82%  "
